$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated analysis timestamp string (shared string in A2)
$ws.Range("A2").Value = "2025-05-22 16:05:34"

# Updated metrics in row 2
$ws.Range("AJ2").Value = 863
$ws.Range("AL2").Value = 17.29128014842301
$ws.Range("AM2").Value = 32.02226345083488
$ws.Range("AN2").Value = 50.68645640074212
$ws.Range("AO2").Value = 1380733.44
$ws.Range("AP2").Value = 259056.88
$ws.Range("AQ2").Value = 86350.56999999999
$ws.Range("AR2").Value = 79.98961428924844
$ws.Range("AS2").Value = 15.00786416107668
$ws.Range("AT2").Value = 5.002521549674893
$ws.Range("AU2").Value = 45.68672617333717
$ws.Range("AV2").Value = 234.2342512259525
$ws.Range("AW2").Value = 559.3112152902993
